# Add data for 2022-03-06:
# - rename sheet / update header label from "through February 25" to "through February 26"
# - update several neighborhood counts in the "February 2022" column (B) and in a few
#   historical columns where new/updated counts appear.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet and update the matching header/title text.
$ws.Name = "Through 2022-02-26"
$ws.Range("B1").Value = "February 2022 (through February 26)"

# Updated counts (row -> neighborhood, for reference):
#   2  Englewood           H2: 1 -> 2
#   3  Austin              F3: 8 -> 9
#   7  Auburn Gresham      L7: 1 -> 2
#   8  North Lawndale      B8: 7 -> 8
#  13  Little Italy, UIC   B13: 2 -> 3
#  16  Chicago Lawn        J16: 1 -> 2 ; L16: 1 -> 2
#  20  Kenwood             B20: 4 -> 5
#  22  Humboldt Park       B22: 2 -> 3
#  23  Washington Heights  D23: (blank) -> 1 ; H23: (blank) -> 1
#  27  Roseland            L27: 1 -> 2
#  33  West Lawn           B33: 2 -> 3
#  34  West Loop           L34: 1 -> 2
#  48  Lake View           B48: 4 -> 5
#  61  Chinatown           B61: 9 -> 11
#  63  Douglas             N63: 1 -> 2

$ws.Range("H2").Value = 2
$ws.Range("F3").Value = 9
$ws.Range("L7").Value = 2
$ws.Range("B8").Value = 8
$ws.Range("B13").Value = 3
$ws.Range("J16").Value = 2
$ws.Range("L16").Value = 2
$ws.Range("B20").Value = 5
$ws.Range("B22").Value = 3
$ws.Range("D23").Value = 1
$ws.Range("H23").Value = 1
$ws.Range("L27").Value = 2
$ws.Range("B33").Value = 3
$ws.Range("L34").Value = 2
$ws.Range("B48").Value = 5
$ws.Range("B61").Value = 11
$ws.Range("N63").Value = 2
